$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALNY")

# Row 17 (EBIT Margin)
$ws.Range("B17").Value = -1.4091

# Row 20 (Free Cash Flow Margin)
$ws.Range("B20").Value = -1.1207
$ws.Range("D20").Value = -1.8555
$ws.Range("E20").Value = -2.4712
$ws.Range("F20").Value = -1.5485
$ws.Range("G20").Value = -1.9048

# Row 28 (EBITDA Margin)
$ws.Range("B28").Value = -1.3215
$ws.Range("D28").Value = -2.2274
$ws.Range("E28").Value = -2.5479
$ws.Range("F28").Value = -3.1928
$ws.Range("G28").Value = -4.0276

# Row 29 (Operating Cash Flow Margin)
$ws.Range("B29").Value = -1.0017
$ws.Range("D29").Value = -1.6373
$ws.Range("E29").Value = -2.1492
$ws.Range("F29").Value = -1.1438
$ws.Range("G29").Value = -1.267
